$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 167
$ws.Range("F5").Value = 3932
$ws.Range("F7").Value = 2591
$ws.Range("F9").Value = 3207
$ws.Range("F11").Value = 2337
$ws.Range("F12").Value = 49
$ws.Range("F15").Value = 468
$ws.Range("F16").Value = 21
$ws.Range("F18").Value = 224
$ws.Range("F19").Value = 357
$ws.Range("F20").Value = 315
$ws.Range("F21").Value = 447
$ws.Range("F22").Value = 681
$ws.Range("F23").Value = 1427
$ws.Range("F24").Value = 65
$ws.Range("F26").Value = 1311
$ws.Range("F27").Value = 142
$ws.Range("F28").Value = 148
$ws.Range("F29").Value = 6
$ws.Range("F30").Value = 73
$ws.Range("F31").Value = 4420
$ws.Range("F32").Value = 4280
$ws.Range("F33").Value = 90
$ws.Range("F34").Value = 211
$ws.Range("F35").Value = 74
$ws.Range("F36").Value = 17
$ws.Range("F37").Value = 1161
$ws.Range("F38").Value = 162
$ws.Range("F40").Value = 497
$ws.Range("F41").Value = 14
$ws.Range("F42").Value = 1325
$ws.Range("F43").Value = 185
$ws.Range("F45").Value = 114
$ws.Range("F46").Value = 44
$ws.Range("F49").Value = 4

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = "不可售"
$ws.Range("F16").Value = 217

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2336
$ws.Range("F5").Value = 14

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 167
$ws.Range("F6").Value = 3932
$ws.Range("F7").Value = 2591
$ws.Range("F9").Value = 3207
$ws.Range("F12").Value = 2337
$ws.Range("F13").Value = 49
$ws.Range("F16").Value = 468
$ws.Range("F17").Value = 21
$ws.Range("F18").Value = 224
$ws.Range("F20").Value = 357
$ws.Range("F21").Value = 681
$ws.Range("F22").Value = 1427
$ws.Range("F23").Value = 1311
$ws.Range("F24").Value = 142
$ws.Range("F26").Value = 73
$ws.Range("F29").Value = 4420
$ws.Range("F30").Value = 4280
$ws.Range("F31").Value = 90
$ws.Range("F32").Value = 17
$ws.Range("F33").Value = 1161
$ws.Range("F34").Value = 162
$ws.Range("F38").Value = 497
$ws.Range("F40").Value = 14
$ws.Range("F43").Value = 1325
$ws.Range("F44").Value = 185
$ws.Range("F45").Value = 114
$ws.Range("F46").Value = 44
$ws.Range("F49").Value = 217
